{"js": "// Add the \"Comments: Sentiment Analysis :\" sub-outline (with its two\n// \"Summary\"/\"Keywords\" sub-bullets) plus two trailing blank list paragraphs\n// right after the existing \"Location Data from status updates\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document starts with exactly one paragraph: \"Location Data from\n// status updates\" - a ListParagraph at ilvl 0 of numId 1. Anchor the new\n// content after it and reuse its list so the new bullets share the same\n// <w:numId>.\nconst anchorParagraph = paragraphs.items[paragraphs.items.length - 1];\nanchorParagraph.load(\"list\");\nawait context.sync();\n\nconst list = anchorParagraph.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\n// 1) \"Comments: Sentiment Analysis :\" - same level (0) as the anchor bullet.\nconst commentsPara = anchorParagraph.insertParagraph(\n  \"Comments: Sentiment Analysis :\",\n  \"After\"\n);\ncommentsPara.style = \"ListParagraph\";\nawait context.sync();\ncommentsPara.attachToList(listId, 0);\nawait context.sync();\n\n// 2) \"Summary\" - nested one level deeper (1).\nconst summaryPara = commentsPara.insertParagraph(\"Summary\", \"After\");\nsummaryPara.style = \"ListParagraph\";\nawait context.sync();\nsummaryPara.attachToList(listId, 1);\nawait context.sync();\n\n// 3) \"Keywords\" - same nested level (1) as \"Summary\".\nconst keywordsPara = summaryPara.insertParagraph(\"Keywords\", \"After\");\nkeywordsPara.style = \"ListParagraph\";\nawait context.sync();\nkeywordsPara.attachToList(listId, 1);\nawait context.sync();\n\n// 4) Empty bullet back at level 0 (still part of the list).\nconst blankListPara = keywordsPara.insertParagraph(\"\", \"After\");\nblankListPara.style = \"ListParagraph\";\nawait context.sync();\nblankListPara.attachToList(listId, 0);\nawait context.sync();\n\n// 5) Trailing empty ListParagraph with no numbering at all.\nconst blankPara = blankListPara.insertParagraph(\"\", \"After\");\nblankPara.style = \"ListParagraph\";\nawait context.sync();\n", "ps1": "# Add the \"Comments: Sentiment Analysis :\" sub-outline (with its two\n# \"Summary\"/\"Keywords\" sub-bullets) plus two trailing blank list paragraphs\n# right after the existing \"Location Data from status updates\" bullet.\n\n$d = $word.ActiveDocument\n\n# The document starts with exactly one paragraph: \"Location Data from\n# status updates\" - a ListParagraph at level 1 (ilvl 0) of the document's\n# only list (numId 1). Anchor the new content after it.\n$anchor = $d.Paragraphs.Last\n\n# 1) \"Comments: Sentiment Analysis :\" - same level (1 / ilvl 0) as the anchor.\n$anchor.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"Comments: Sentiment Analysis :\"\n$p1.Style = \"ListParagraph\"\n$p1.Range.ListFormat.ListLevelNumber = 1\n\n# 2) \"Summary\" - nested one level deeper (2 / ilvl 1).\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"Summary\"\n$p2.Style = \"ListParagraph\"\n$p2.Range.ListFormat.ListLevelNumber = 2\n\n# 3) \"Keywords\" - same nested level (2 / ilvl 1) as \"Summary\".\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Last\n$p3.Range.Text = \"Keywords\"\n$p3.Style = \"ListParagraph\"\n$p3.Range.ListFormat.ListLevelNumber = 2\n\n# 4) Empty bullet back at level 1 (ilvl 0), still part of the list.\n$p3.Range.InsertParagraphAfter()\n$p4 = $d.Paragraphs.Last\n$p4.Style = \"ListParagraph\"\n$p4.Range.ListFormat.ListLevelNumber = 1\n\n# 5) Trailing empty ListParagraph with no numbering at all.\n$p4.Range.InsertParagraphAfter()\n$p5 = $d.Paragraphs.Last\n$p5.Style = \"ListParagraph\"\n"}
